# Applies the Coinranking "cryptos" price/volume refresh captured in the
# GitHub Actions commit of Wed Jan 31 11:22:47 UTC 2024.
# Only column D (Price) and E (Volume(1h)) values change for most rows;
# rows 44-45 additionally swap which coin (VeChain / EnergySwap) occupies
# that rank, including columns B (Coin) and C (Link).
#
# Price-column values that are valid bare numbers (e.g. "34.30", "1.00")
# get a leading apostrophe so Excel keeps them as literal text instead of
# normalising them to a Number and dropping the trailing zero - exactly
# like the "Qualify as text" behaviour a user gets by typing 'value in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.485.36'
$ws.Range("E2").Value = '  -2.16%  '

# Row 3
$ws.Range("E3").Value = '  -1.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''302.12'
$ws.Range("E5").Value = '  -2.68%  '

# Row 6
$ws.Range("D6").Value = '''98.13'
$ws.Range("E6").Value = '  -6.03%  '

# Row 7
$ws.Range("D7").Value = '''0.505'
$ws.Range("E7").Value = '  -5.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").Value = '''0.498'
$ws.Range("E9").Value = '  -5.84%  '

# Row 10
$ws.Range("D10").Value = '''34.30'
$ws.Range("E10").Value = '  -6.64%  '

# Row 11
$ws.Range("E11").Value = '  -3.20%  '

# Row 12
$ws.Range("D12").Value = '''0.112'
$ws.Range("E12").Value = '  -0.01%  '

# Row 13
$ws.Range("E13").Value = '  -4.42%  '

# Row 14
$ws.Range("D14").Value = '2.644.30'
$ws.Range("E14").Value = '  -1.05%  '

# Row 15
$ws.Range("D15").Value = '''15.58'
$ws.Range("E15").Value = '  +2.82%  '

# Row 16
$ws.Range("D16").Value = '2.294.89'
$ws.Range("E16").Value = '  -0.83%  '

# Row 17
$ws.Range("E17").Value = '  -1.23%  '

# Row 18
$ws.Range("D18").Value = '42.421.98'
$ws.Range("E18").Value = '  -2.12%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0898'
$ws.Range("E19").Value = '  -2.97%  '

# Row 20
$ws.Range("D20").Value = '''11.45'
$ws.Range("E20").Value = '  -5.92%  '

# Row 21
$ws.Range("E21").Value = '  -2.40%  '

# Row 22
$ws.Range("D22").Value = '''67.68'
$ws.Range("E22").Value = '  -0.83%  '

# Row 23
$ws.Range("D23").Value = '''234.48'
$ws.Range("E23").Value = '  -3.24%  '

# Row 24
$ws.Range("D24").Value = '''1.96'
$ws.Range("E24").Value = '  -3.48%  '

# Row 25
$ws.Range("D25").Value = '''2.52'
$ws.Range("E25").Value = '  -3.63%  '

# Row 26
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.09%  '

# Row 27
$ws.Range("D27").Value = '''24.91'
$ws.Range("E27").Value = '  +0.19%  '

# Row 28
$ws.Range("D28").Value = '''2.29'
$ws.Range("E28").Value = '  -4.19%  '

# Row 29
$ws.Range("D29").Value = '''34.60'
$ws.Range("E29").Value = '  -6.75%  '

# Row 30
$ws.Range("D30").Value = '''9.15'
$ws.Range("E30").Value = '  -5.11%  '

# Row 31
$ws.Range("D31").Value = '''163.04'
$ws.Range("E31").Value = '  -2.24%  '

# Row 32
$ws.Range("E32").Value = '  +0.05%  '

# Row 33
$ws.Range("D33").Value = '''5.00'
$ws.Range("E33").Value = '  -5.42%  '

# Row 34
$ws.Range("D34").Value = '''4.60'
$ws.Range("E34").Value = '  +1.78%  '

# Row 35
$ws.Range("E35").Value = '  -4.87%  '

# Row 36
$ws.Range("D36").Value = '''0.0711'
$ws.Range("E36").Value = '  -4.47%  '

# Row 37
$ws.Range("D37").Value = '''16.91'
$ws.Range("E37").Value = '  -7.58%  '

# Row 38
$ws.Range("D38").Value = '''2.87'
$ws.Range("E38").Value = '  -6.13%  '

# Row 39
$ws.Range("E39").Value = '  -4.93%  '

# Row 40
$ws.Range("D40").Value = '''0.101'
$ws.Range("E40").Value = '  -5.12%  '

# Row 41
$ws.Range("E41").Value = '  -3.87%  '

# Row 42
$ws.Range("D42").Value = '''2.44'
$ws.Range("E42").Value = '  -10.08%  '

# Row 43
$ws.Range("D43").Value = '1.971.96'
$ws.Range("E43").Value = '  -1.00%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0279'
$ws.Range("E44").Value = '  -4.77%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''18.49'
$ws.Range("E45").Value = '  -2.66%  '

# Row 46
$ws.Range("D46").Value = '''10.13'
$ws.Range("E46").Value = '  +1.16%  '

# Row 47
$ws.Range("D47").Value = '''2.89'
$ws.Range("E47").Value = '  -5.36%  '

# Row 48
$ws.Range("D48").Value = '''55.37'
$ws.Range("E48").Value = '  -2.96%  '

# Row 49
$ws.Range("E49").Value = '  -3.82%  '

# Row 50
$ws.Range("D50").Value = '2.515.80'
$ws.Range("E50").Value = '  -0.90%  '
